# Generate Report for Archive
#
# The localization status for the two handed-back files moved from
# "Ready for handoff" to "In Translation" on the Overview sheet (zh-cn and
# de-de status columns) as well as on each per-locale status sheet. Excel
# then re-flowed the (now narrower) status columns to fit the shorter text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) and de-de (col F) status columns -------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"

# --- Per-locale sheets: Status column (col C) ------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"

# --- Re-fit the affected columns to the new (shorter) text ----------------
# ColumnWidth is expressed in characters; 12.5 is the closest representable
# value to the post-edit autofit width used when the report was archived.
$overview.Columns("E:F").ColumnWidth = 12.5
$zhcn.Columns("C:C").ColumnWidth = 12.5
$dede.Columns("C:C").ColumnWidth = 12.5
